# Applies a re-shuffle of the weekly records (columns D, J, K, L, M, N, O, P, Q)
# among rows 2-15 of the "Hortaliza, Macroferia Regional de Talca - Cilantro" sheet.
# Rows 13 and 15 are unchanged; all other rows receive another row's values
# for the date/volume/price/unit/origin columns, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (target) values for each affected row, for columns D,J,K,L,M,N,O,P,Q
$data = @{
    2  = @{ D = 44357; J = 150; K = 6500; L = 6500; M = 6500; N = '$/caja 20 docenas'; O = 'Región del Maule';     P = 6500; Q = 1  }
    3  = @{ D = 44364; J = 100; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región Metropolitana'; P = 194;  Q = 36 }
    4  = @{ D = 44386; J = 200; K = 6500; L = 6500; M = 6500; N = '$/caja 36 atados';  O = 'Región Metropolitana'; P = 181;  Q = 36 }
    5  = @{ D = 44354; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región del Maule';     P = 194;  Q = 36 }
    6  = @{ D = 44355; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región Metropolitana'; P = 194;  Q = 36 }
    7  = @{ D = 44372; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región Metropolitana'; P = 194;  Q = 36 }
    8  = @{ D = 44340; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región del Maule';     P = 194;  Q = 36 }
    9  = @{ D = 44371; J = 150; K = 6500; L = 6500; M = 6500; N = '$/caja 36 atados';  O = 'Región Metropolitana'; P = 181;  Q = 36 }
    10 = @{ D = 44362; J = 100; K = 6500; L = 6500; M = 6500; N = '$/caja 36 atados';  O = 'Región Metropolitana'; P = 181;  Q = 36 }
    11 = @{ D = 44358; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región Metropolitana'; P = 194;  Q = 36 }
    12 = @{ D = 44348; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región del Maule';     P = 194;  Q = 36 }
    14 = @{ D = 44342; J = 150; K = 7000; L = 7000; M = 7000; N = '$/caja 36 atados';  O = 'Región del Maule';     P = 194;  Q = 36 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
}
